# Word document is closed after reading
# Update a handful of search-value regex strings in the "invoice_sales" sheet
# (they gained a trailing ";1" / ";" marker), and restore the sheet's view
# (top-left cell / selected cell) to its pre-scroll state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invoice_sales")

$ws.Range("G6").Value = '(\d+\s*)+[,]*\d{2,3};;1'
$ws.Range("H6").Value = '(\d+\s*)+[,]*\d{2,3};;1'
$ws.Range("F4").Value = '\bUSD\b;1;1'
$ws.Range("F5").Value = '((\d{1,3}\s)+\d{1,3},\d{2});1;1'
$ws.Range("E4").Value = '([a-z]{3}\/[a-z]{3});;1'
$ws.Range("E5").Value = '^(\d{1,3},\d{3});1;1'
$ws.Range("E6").Value = '^(\d{1,3},\d{3});;1'
$ws.Range("G4").Value = '\b[net]{3}\b;1;1'

# Restore view/selection state on the sheet
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("G5").Select()
